$d = $word.ActiveDocument

# The timeline information is held in the single table in the document.
# Add a new row that re-uses the formatting of the existing last row
# (Word duplicates the preceding row's cell formatting/shading/widths
# automatically when a row is appended this way).
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

$enDash = [char]0x2013

$newRow.Cells.Item(1).Range.Text = "31/12/2021"
$newRow.Cells.Item(2).Range.Text = "4 Hours 15 Minutes"
$newRow.Cells.Item(3).Range.Text = "World Generation " + $enDash + " Objective 2"
$newRow.Cells.Item(4).Range.Text = "Added province rendering to show all the generated province boundaries. Early into this procedure issues were found in which the provinces were displaying as overlapping sets, but this was identified as a rotation issue. Rotating the provinces by 180 degrees x and 180 degrees y created a fully mapped provincial system. As of now, there are still oceans included within the province definitions " + $enDash + " as well as some small provinces that could be put to better use by combining them with other nearby provinces."
